$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held the professor name ("519033 - Carlos Yujiro
# Shigue") in columns B/C is removed; everything below shifts up one row.
$ws.Rows("13").Delete()

# "Objetivos:" (row 10) now answers with the professor name instead of
# the old free-text objectives paragraph.
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# "Programa resumido:" (was row 14, now row 13 after the shift) becomes
# "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "Programa:" (was row 16, now row 15) loses its long free-text program
# and instead shows the activation date (as text, not an Excel date).
$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Rows("15").RowHeight = 120

# "Método:" (was row 19, now row 18) also gets the professor name.
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

# "Critério:" (was row 20, now row 19) takes over the "Aulas expositivas
# ..." method text that used to belong to "Método:".
$ws.Range("B19").Value = "Aulas expositivas, seminários e práticas ministradas em laboratório."
$ws.Range("C19").Value = "Aulas expositivas, seminários e práticas ministradas em laboratório."

# "Norma de recuperação:" (was row 21, now row 20) takes over the
# weighted-average grading criterion text that used to belong to
# "Critério:".
$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

# "Bibliografia:" (was row 22, now row 21) takes over the recovery-exam
# rule text that used to belong to "Norma de recuperação:", dropping the
# old reference list.
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# Row heights of the two rows that grew to the "long text" band.
$ws.Rows("21").RowHeight = 120
